$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$genes = @(
    "ABCG2",
    "ACSL6",
    "ADAM17",
    "AHSP",
    "AQP9",
    "ARG1",
    "BLVRB",
    "BRCA2",
    "BSG",
    "C1R",
    "CD209",
    "CD47",
    "CD9",
    "CFH",
    "CLEC5A",
    "CPT1A",
    "CROT",
    "EGF",
    "EMP1",
    "FECH",
    "FSCN1",
    "FZD5",
    "G0S2",
    "GPX2",
    "GSTM3",
    "HBD",
    "HEBP1",
    "HMBS",
    "ICAM1",
    "INHBB",
    "JAM3",
    "MAFF",
    "MAP3K6",
    "MARCH8",
    "ME1",
    "MEIS1",
    "MYL4",
    "MYL9",
    "NAP1L2",
    "NEK2",
    "PEX7",
    "PTP4A3",
    "RAD51C",
    "RAP1GAP",
    "RBM38",
    "RETN",
    "RPL14",
    "SDC2",
    "SDHD",
    "SGCD",
    "SLC22A4",
    "SPTB",
    "SSR3",
    "TIMP1",
    "TMEM176A",
    "TMEM176B",
    "TNS1",
    "TPM1",
    "TRIM10",
    "TRIM25",
    "TSC22D1",
    "TSPAN4",
    "VNN1",
    "DACT1",
    "DBNDD1",
    "ERC1",
    "GPR25",
    "NLRP2",
    "RNASE6",
    "SLC14A1",
    "TTTY15",
    "AGAP1",
    "ALPL",
    "APOBEC3G",
    "BTN3A2",
    "CD79B",
    "CFD",
    "CHI3L1",
    "ECHDC3",
    "ERO1B",
    "FAM13A",
    "GALNT6",
    "GPM6A",
    "GPX1",
    "HLA-DPA1",
    "HNRNPH1",
    "IL12RB2",
    "ITPR1",
    "KLHDC4",
    "LGALSL",
    "MRC2",
    "NEBL",
    "NSUN6",
    "PCGF3",
    "POLR1D",
    "PRKAR2B",
    "RUNX3",
    "SLC7A8",
    "SNORA21",
    "SPATA20",
    "STRN3",
    "TPP2",
    "YME1L1",
    "ACAA2",
    "ACP1",
    "ACP2",
    "ADM",
    "ALB",
    "ASNS",
    "BPGM",
    "CACNA2D2",
    "CALD1",
    "CCND1",
    "CD28",
    "CDC25B",
    "DGAT1",
    "DHFR",
    "DHX16",
    "DLG5",
    "EBP",
    "ENDOD1",
    "ENO3",
    "ETV1",
    "FDXR",
    "GATA1",
    "GDE1",
    "GGA2",
    "GNG11",
    "GYPA",
    "GYPB",
    "HBQ1",
    "HDGF",
    "HIST1H1T",
    "HLA-F",
    "IRF5",
    "KLHDC3",
    "LAMP2",
    "MAP7",
    "MARK3",
    "MCM5",
    "MKRN1",
    "MMP1",
    "MPP1",
    "MPZL2",
    "OPN3",
    "PCDH7",
    "PDK3",
    "PGRMC1",
    "PLOD2",
    "PMPCA",
    "PPP2R5B",
    "PQLC1",
    "PREB",
    "RIOK3",
    "SEC14L1",
    "SLC25A37",
    "SMTN",
    "SNCA",
    "SNRPA",
    "THBD",
    "TMCC2",
    "TMEM158",
    "TRIM58",
    "TSPAN5",
    "XYLT2",
    "BEX3",
    "CLEC10A",
    "CNTN5",
    "ENY2",
    "FHL2",
    "LAMTOR2",
    "DOC2A",
    "MARC1",
    "PASK",
    "PLCL1",
    "PRR5L",
    "TNFRSF21",
    "ACADVL",
    "CSF1",
    "DUSP1",
    "DUT",
    "IDH1",
    "ITGB1",
    "MCL1",
    "NFE2",
    "PGK1",
    "PNP",
    "RB1",
    "RXRA",
    "SLC6A8",
    "SOCS1",
    "SPTBN2",
    "TAL1",
    "TCF7L2",
    "TFDP2",
    "TGIF2",
    "TNPO2",
    "UQCR10",
    "VCAN",
    "VEGFA",
    "ATP2A3",
    "GZMH",
    "MAPKAPK2",
    "RGS9",
    "ZSWIM8"
)

for ($i = 0; $i -lt $genes.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $genes[$i]
}
